$d = $word.ActiveDocument

# --- locate the target paragraph -------------------------------------------------
# The doc has several "PRIMARY KEY('student_id', 'telefoonnummer')," paragraphs;
# the one we want is the one that is shortly followed (4 paragraphs later) by the
# "DROP TABLE IF EXISTS ... 'semester'" paragraph that currently owns the
# _GoBack bookmark.
$targetIndex = 0
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "PRIMARY KEY(‘student_id’, ‘telefoonnummer’),*") {
        $t2 = $d.Paragraphs.Item($i + 4).Range.Text
        if ($t2 -like "DROP TABLE IF EXISTS*semester*") {
            $targetIndex = $i
        }
    }
}

$para = $d.Paragraphs.Item($targetIndex)

# --- replace 'telefoonnummer' with 'email' (quotes included so the run merges
#     and the proofErr spell-check markup around the word is cleared) -----------
$searchRange = $para.Range.Duplicate
$searchRange.Find.Execute("‘telefoonnummer’", $true, $false, $false, $false, $false, $true, 1, $false, "‘email’", 2) | Out-Null

# --- find the new "email" run so we know exactly where to split/bookmark ------
$emailRange = $para.Range.Duplicate
$emailRange.Find.Execute("email", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$emailStart = $emailRange.Start
$emailEnd = $emailRange.End

# --- split "...', 'email" into separate runs by dropping a temporary bookmark
#     right before "email" (bookmarks force a run boundary in this engine) -----
$d.Bookmarks.Add("ZZTmpSplit", $d.Range($emailStart, $emailStart)) | Out-Null

# --- move the _GoBack bookmark from the old "DROP TABLE IF EXISTS" location to
#     right after the new "email" run ------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($emailEnd, $emailEnd)) | Out-Null

# --- clean up the temporary bookmark (its run split remains in place) ---------
$d.Bookmarks.Item("ZZTmpSplit").Delete()
